{"js": "const replacements = [\n  [\"361\u00f75=\", \"750\u00f73=\"],\n  [\"391\u00f74=\", \"815\u00f76=\"],\n  [\"791\u00f79=\", \"952\u00f73=\"],\n  [\"261\u00f76=\", \"611\u00f76=\"],\n  [\"408\u00f77=\", \"184\u00f72=\"],\n  [\"263\u00f73=\", \"369\u00f74=\"],\n  [\"344\u00f79=\", \"814\u00f75=\"],\n  [\"956\u00f78=\", \"542\u00f77=\"],\n  [\"183\u00f72=\", \"745\u00f75=\"],\n  [\"697\u00f75=\", \"260\u00f74=\"],\n  [\"660\u00f74=\", \"771\u00f74=\"],\n  [\"807\u00f76=\", \"682\u00f72=\"],\n  [\"160\u00f75=\", \"140\u00f75=\"],\n  [\"253\u00f77=\", \"553\u00f77=\"],\n  [\"175\u00f79=\", \"322\u00f77=\"],\n  [\"315\u00f74=\", \"111\u00f79=\"],\n  [\"502\u00f78=\", \"715\u00f73=\"],\n  [\"923\u00f76=\", \"388\u00f72=\"],\n  [\"519\u00f75=\", \"267\u00f79=\"],\n  [\"921\u00f78=\", \"379\u00f76=\"],\n  [\"584\u00f76=\", \"861\u00f72=\"],\n  [\"454\u00f72=\", \"991\u00f76=\"],\n  [\"121\u00f79=\", \"991\u00f79=\"],\n  [\"120\u00f78=\", \"186\u00f72=\"],\n  [\"464\u00f75=\", \"387\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"361\u00f75=\"; new=\"750\u00f73=\"},\n    @{old=\"391\u00f74=\"; new=\"815\u00f76=\"},\n    @{old=\"791\u00f79=\"; new=\"952\u00f73=\"},\n    @{old=\"261\u00f76=\"; new=\"611\u00f76=\"},\n    @{old=\"408\u00f77=\"; new=\"184\u00f72=\"},\n    @{old=\"263\u00f73=\"; new=\"369\u00f74=\"},\n    @{old=\"344\u00f79=\"; new=\"814\u00f75=\"},\n    @{old=\"956\u00f78=\"; new=\"542\u00f77=\"},\n    @{old=\"183\u00f72=\"; new=\"745\u00f75=\"},\n    @{old=\"697\u00f75=\"; new=\"260\u00f74=\"},\n    @{old=\"660\u00f74=\"; new=\"771\u00f74=\"},\n    @{old=\"807\u00f76=\"; new=\"682\u00f72=\"},\n    @{old=\"160\u00f75=\"; new=\"140\u00f75=\"},\n    @{old=\"253\u00f77=\"; new=\"553\u00f77=\"},\n    @{old=\"175\u00f79=\"; new=\"322\u00f77=\"},\n    @{old=\"315\u00f74=\"; new=\"111\u00f79=\"},\n    @{old=\"502\u00f78=\"; new=\"715\u00f73=\"},\n    @{old=\"923\u00f76=\"; new=\"388\u00f72=\"},\n    @{old=\"519\u00f75=\"; new=\"267\u00f79=\"},\n    @{old=\"921\u00f78=\"; new=\"379\u00f76=\"},\n    @{old=\"584\u00f76=\"; new=\"861\u00f72=\"},\n    @{old=\"454\u00f72=\"; new=\"991\u00f76=\"},\n    @{old=\"121\u00f79=\"; new=\"991\u00f79=\"},\n    @{old=\"120\u00f78=\"; new=\"186\u00f72=\"},\n    @{old=\"464\u00f75=\"; new=\"387\u00f73=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    [void]$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $r.new, 2)\n}\n"}
